$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 683.3333
$ws.Range("I6").Value = 675
$ws.Range("K6").Value = 2025
$ws.Range("M6").Value = -1913
$ws.Range("H40").Value = 1664.7059
$ws.Range("I40").Value = 1555
$ws.Range("J40").Value = 1821.4286
$ws.Range("K40").Value = 1555
$ws.Range("L40").Value = 1821.4286
$ws.Range("M40").Value = -1380
$ws.Range("N40").Value = -2171.4286
$ws.Range("H64").Value = 3083.3333
$ws.Range("I64").Value = 2918.182
$ws.Range("J64").Value = 3223.077
$ws.Range("K64").Value = 2918.182
$ws.Range("L64").Value = 3223.077
$ws.Range("M64").Value = -2670.182
$ws.Range("N64").Value = -3719.077
$ws.Range("H67").Value = 3083.3333
$ws.Range("I67").Value = 2918.182
$ws.Range("J67").Value = 3223.077
$ws.Range("K67").Value = 2918.182
$ws.Range("L67").Value = 3223.077
$ws.Range("M67").Value = -2060.182
$ws.Range("N67").Value = -4939.077
$ws.Range("H76").Value = 7977.4517
$ws.Range("I76").Value = 13708.417
$ws.Range("J76").Value = 4357.8945
$ws.Range("K76").Value = 13708.417
$ws.Range("L76").Value = 4357.8945
$ws.Range("M76").Value = -13393.417
$ws.Range("N76").Value = -4987.8945
$ws.Range("H79").Value = 7977.4517
$ws.Range("I79").Value = 13708.417
$ws.Range("J79").Value = 4357.8945
$ws.Range("K79").Value = 13708.417
$ws.Range("L79").Value = 4357.8945
$ws.Range("M79").Value = -12616.417
$ws.Range("N79").Value = -6541.8945
$ws.Range("H132").Value = 2226
$ws.Range("I132").Value = 1630.8413
$ws.Range("J132").Value = 4569.4375
$ws.Range("K132").Value = 4892.5239
$ws.Range("L132").Value = 13708.3125
$ws.Range("M132").Value = -2362.5239
$ws.Range("N132").Value = -18768.3125

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 30799.334
$ws.Range("J44").Value = 30799.334
$ws.Range("L44").Value = 30799.334
$ws.Range("N44").Value = -31775.334
$ws.Range("H55").Value = 30448.5
$ws.Range("J55").Value = 39315.332
$ws.Range("L55").Value = 39315.332
$ws.Range("N55").Value = -39945.332
$ws.Range("H63").Value = 4166.5835
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 4555.4443
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 4555.4443
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -5927.4443
$ws.Range("H66").Value = 4166.5835
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 4555.4443
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 22777.2215
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -29641.2215
$ws.Range("H80").Value = 32264
$ws.Range("J80").Value = 38055
$ws.Range("L80").Value = 38055
$ws.Range("N80").Value = -40051
$ws.Range("H83").Value = 32264
$ws.Range("J83").Value = 38055
$ws.Range("L83").Value = 114165
$ws.Range("N83").Value = -124149
$ws.Range("H109").Value = 28059.5
$ws.Range("J109").Value = 28059.5
$ws.Range("L109").Value = 28059.5
$ws.Range("N109").Value = -30833.5
$ws.Range("H122").Value = 13496.588
$ws.Range("I122").Value = 21529.4
$ws.Range("K122").Value = 64588.2
$ws.Range("M122").Value = -62138.2

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 31624.666
$ws.Range("J35").Value = 31624.666
$ws.Range("L35").Value = 31624.666
$ws.Range("N35").Value = -32244.666
$ws.Range("H82").Value = 34254.668
$ws.Range("J82").Value = 35307.6
$ws.Range("L82").Value = 35307.6
$ws.Range("N82").Value = -36073.6
$ws.Range("H85").Value = 34254.668
$ws.Range("J85").Value = 35307.6
$ws.Range("L85").Value = 35307.6
$ws.Range("N85").Value = -37959.6
$ws.Range("H86").Value = 1713.8182
$ws.Range("I86").Value = 1673.0952
$ws.Range("J86").Value = 1785.0834
$ws.Range("K86").Value = 1673.0952
$ws.Range("L86").Value = 1785.0834
$ws.Range("M86").Value = -550.0952
$ws.Range("N86").Value = -4031.0834
$ws.Range("H89").Value = 1713.8182
$ws.Range("I89").Value = 1673.0952
$ws.Range("J89").Value = 1785.0834
$ws.Range("K89").Value = 8365.476000000001
$ws.Range("L89").Value = 8925.416999999999
$ws.Range("M89").Value = -2749.476000000001
$ws.Range("N89").Value = -20157.417
$ws.Range("H105").Value = 991632.4
$ws.Range("I105").Value = 1750447.9
$ws.Range("J105").Value = 5172.2
$ws.Range("K105").Value = 1750447.9
$ws.Range("L105").Value = 5172.2
$ws.Range("M105").Value = -1748700.9
$ws.Range("N105").Value = -8666.200000000001
$ws.Range("H108").Value = 38866
$ws.Range("J108").Value = 38866
$ws.Range("L108").Value = 38866
$ws.Range("N108").Value = -46546

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 19915
$ws.Range("I41").Value = 8000
$ws.Range("J41").Value = 23886.666
$ws.Range("K41").Value = 8000
$ws.Range("L41").Value = 23886.666
$ws.Range("M41").Value = -7572
$ws.Range("N41").Value = -24742.666
$ws.Range("H60").Value = 24815.928
$ws.Range("J60").Value = 24815.928
$ws.Range("L60").Value = 24815.928
$ws.Range("N60").Value = -25837.928
$ws.Range("H62").Value = 3749.6667
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 3749.6667
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H109").Value = 13852.223
$ws.Range("J109").Value = 13852.223
$ws.Range("L109").Value = 13852.223
$ws.Range("N109").Value = -15932.223
$ws.Range("H132").Value = 2150.121
$ws.Range("I132").Value = 1742.7778
$ws.Range("J132").Value = 3983.1667
$ws.Range("K132").Value = 5228.3334
$ws.Range("L132").Value = 11949.5001
$ws.Range("M132").Value = -2698.3334
$ws.Range("N132").Value = -17009.5001

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 44384.19
$ws.Range("I70").Value = 56265
$ws.Range("K70").Value = 56265
$ws.Range("M70").Value = -55995
$ws.Range("H73").Value = 44384.19
$ws.Range("I73").Value = 56265
$ws.Range("K73").Value = 56265
$ws.Range("M73").Value = -55329
$ws.Range("H80").Value = 2237.5
$ws.Range("I80").Value = 2333.3333
$ws.Range("J80").Value = 1950
$ws.Range("K80").Value = 2333.3333
$ws.Range("L80").Value = 1950
$ws.Range("M80").Value = -1335.3333
$ws.Range("N80").Value = -3946
$ws.Range("H83").Value = 2237.5
$ws.Range("I83").Value = 2333.3333
$ws.Range("J83").Value = 1950
$ws.Range("K83").Value = 11666.6665
$ws.Range("L83").Value = 9750
$ws.Range("M83").Value = -6674.666499999999
$ws.Range("N83").Value = -19734
$ws.Range("H122").Value = 1856.7273
$ws.Range("I122").Value = 1840.9231
$ws.Range("J122").Value = 1915.4286
$ws.Range("K122").Value = 5522.7693
$ws.Range("L122").Value = 5746.2858
$ws.Range("M122").Value = -3072.7693
$ws.Range("N122").Value = -10646.2858
$ws.Range("H123").Value = 32973.25
$ws.Range("J123").Value = 32973.25
$ws.Range("L123").Value = 32973.25
$ws.Range("N123").Value = -37873.25
$ws.Range("H132").Value = 2459.8572
$ws.Range("I132").Value = 2314.8333
$ws.Range("K132").Value = 6944.499899999999
$ws.Range("M132").Value = -4414.499899999999

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3995.75
$ws.Range("I132").Value = 4199.1665
$ws.Range("J132").Value = 3385.5
$ws.Range("K132").Value = 12597.4995
$ws.Range("L132").Value = 10156.5
$ws.Range("M132").Value = -10067.4995
$ws.Range("N132").Value = -15216.5
$ws.Range("H136").Value = 2780.6
$ws.Range("I136").Value = 2659.7727
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 7979.3181
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -5429.3181
$ws.Range("N136").Value = -16100.0001

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 18188.5
$ws.Range("J109").Value = 18188.5
$ws.Range("L109").Value = 18188.5
$ws.Range("N109").Value = -20962.5
$ws.Range("H122").Value = 12466.619
$ws.Range("I122").Value = 20300.363
$ws.Range("J122").Value = 3849.5
$ws.Range("K122").Value = 60901.08900000001
$ws.Range("L122").Value = 11548.5
$ws.Range("M122").Value = -58451.08900000001
$ws.Range("N122").Value = -16448.5
